$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run of SGNN dialog-act annotation after transcript clean up.
# Each entry: spreadsheet row number, new DAMSLTag (col I), new DialogAct label (col J)
$updates = @(
    @{ Row = 3; I = "%"; J = "Uninterpretable" },
    @{ Row = 4; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 6; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 19; I = "b"; J = "Acknowledge (Backchannel)" },
    @{ Row = 34; I = "sv"; J = "Statement-opinion" },
    @{ Row = 35; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 53; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 58; I = "ba"; J = "Appreciation" },
    @{ Row = 67; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 76; I = "sv"; J = "Statement-opinion" },
    @{ Row = 83; I = "aa"; J = "Agree/Accept" },
    @{ Row = 85; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 96; I = "aa"; J = "Agree/Accept" },
    @{ Row = 101; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 102; I = "sv"; J = "Statement-opinion" },
    @{ Row = 106; I = "%"; J = "Uninterpretable" },
    @{ Row = 112; I = "%"; J = "Uninterpretable" },
    @{ Row = 115; I = "aa"; J = "Agree/Accept" },
    @{ Row = 122; I = "sv"; J = "Statement-opinion" },
    @{ Row = 125; I = "aa"; J = "Agree/Accept" },
    @{ Row = 128; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 137; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 139; I = "sv"; J = "Statement-opinion" },
    @{ Row = 143; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 144; I = "b"; J = "Acknowledge (Backchannel)" },
    @{ Row = 150; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 154; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 156; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 169; I = "aa"; J = "Agree/Accept" },
    @{ Row = 175; I = "ba"; J = "Appreciation" },
    @{ Row = 178; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 196; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 197; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 198; I = "sv"; J = "Statement-opinion" },
    @{ Row = 244; I = "sv"; J = "Statement-opinion" },
    @{ Row = 248; I = "sv"; J = "Statement-opinion" },
    @{ Row = 252; I = "%"; J = "Uninterpretable" },
    @{ Row = 263; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 277; I = "aa"; J = "Agree/Accept" },
    @{ Row = 295; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 299; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 305; I = "sv"; J = "Statement-opinion" },
    @{ Row = 324; I = "aa"; J = "Agree/Accept" },
    @{ Row = 325; I = "aa"; J = "Agree/Accept" },
    @{ Row = 328; I = "aa"; J = "Agree/Accept" },
    @{ Row = 332; I = "aa"; J = "Agree/Accept" },
    @{ Row = 333; I = "ba"; J = "Appreciation" },
    @{ Row = 340; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 366; I = "b"; J = "Acknowledge (Backchannel)" },
    @{ Row = 367; I = "sv"; J = "Statement-opinion" },
    @{ Row = 368; I = "%"; J = "Uninterpretable" },
    @{ Row = 371; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 373; I = "b"; J = "Acknowledge (Backchannel)" },
    @{ Row = 394; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 409; I = "sv"; J = "Statement-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

Write-Host "Applied $($updates.Count) dialog act updates."
